# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect newly refreshed totals, as generated by the gh-pages build.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows indexed by event)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 46
$wsExhibit.Range("F4").Value = 3471
$wsExhibit.Range("F5").Value = 2168
$wsExhibit.Range("F6").Value = 422
$wsExhibit.Range("F8").Value = 56
$wsExhibit.Range("F9").Value = 46
$wsExhibit.Range("F10").Value = 1287
$wsExhibit.Range("F12").Value = 1677
$wsExhibit.Range("F13").Value = 120

# Sheet "全部类型" (same events, different row offsets)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 46
$wsAll.Range("F4").Value = 3471
$wsAll.Range("F5").Value = 2168
$wsAll.Range("F6").Value = 422
$wsAll.Range("F9").Value = 56
$wsAll.Range("F10").Value = 46
$wsAll.Range("F13").Value = 1287
$wsAll.Range("F15").Value = 1677
$wsAll.Range("F16").Value = 120
